$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 197
$rng = $ws.Range("A197:H197")
$rng.NumberFormat = "@"
$ws.Range("A197").Value = "WC48 P5F"
$ws.Range("B197").Value = "Etiquetadora"
$ws.Range("C197").Value = "2024-06-11"
$ws.Range("D197").Value = "09:18:43"
$ws.Range("E197").Value = "Mañana"
$ws.Range("F197").Value = "09:18:44"
$ws.Range("G197").Value = "0:00:01"
$ws.Range("H197").Value = "-0.01 minutos"
$rng.ClearFormats()

# Row 198
$rng = $ws.Range("A198:H198")
$rng.NumberFormat = "@"
$ws.Range("A198").Value = "WC48 P5F"
$ws.Range("B198").Value = "Fallo etiqueta"
$ws.Range("C198").Value = "2024-06-11"
$ws.Range("D198").Value = "09:18:49"
$ws.Range("E198").Value = "Mañana"
$ws.Range("F198").Value = "09:18:50"
$ws.Range("G198").Value = "0:00:01"
$ws.Range("H198").Value = "0.02 minutos"
$rng.ClearFormats()

# Row 199
$rng = $ws.Range("A199:H199")
$rng.NumberFormat = "@"
$ws.Range("A199").Value = "WC48 P5F"
$ws.Range("B199").Value = "Detección de sealling mal puesto"
$ws.Range("C199").Value = "2024-06-11"
$ws.Range("D199").Value = "09:19:06"
$ws.Range("E199").Value = "Mañana"
$ws.Range("F199").Value = "09:19:07"
$ws.Range("G199").Value = "0:00:01"
$ws.Range("H199").Value = "0.06 minutos"
$rng.ClearFormats()

# Row 200
$rng = $ws.Range("A200:H200")
$rng.NumberFormat = "@"
$ws.Range("A200").Value = "WC49 P5H"
$ws.Range("B200").Value = "Screw K30 no lo detecta puesto"
$ws.Range("C200").Value = "2024-06-11"
$ws.Range("D200").Value = "09:28:36"
$ws.Range("E200").Value = "Mañana"
$ws.Range("F200").Value = "09:28:38"
$ws.Range("G200").Value = "0:00:02"
$ws.Range("H200").Value = "-0.00 minutos"
$rng.ClearFormats()

# Row 201
$rng = $ws.Range("A201:H201")
$rng.NumberFormat = "@"
$ws.Range("A201").Value = "WC49 P5H"
$ws.Range("B201").Value = "Power atascado en prensa, cuesta sacar"
$ws.Range("C201").Value = "2024-06-11"
$ws.Range("D201").Value = "09:29:11"
$ws.Range("E201").Value = "Mañana"
$ws.Range("F201").Value = "09:29:11"
$ws.Range("G201").Value = "0:00:00"
$ws.Range("H201").Value = "0.14 minutos"
$rng.ClearFormats()

# Row 202
$rng = $ws.Range("A202:H202")
$rng.NumberFormat = "@"
$ws.Range("A202").Value = "WC49 P5H"
$ws.Range("B202").Value = "Etiquetadora"
$ws.Range("C202").Value = "2024-06-11"
$ws.Range("D202").Value = "09:29:14"
$ws.Range("E202").Value = "Mañana"
$ws.Range("F202").Value = "09:29:14"
$ws.Range("G202").Value = "0:00:00"
$ws.Range("H202").Value = "0.10 minutos"
$rng.ClearFormats()

# Row 203
$rng = $ws.Range("A203:H203")
$rng.NumberFormat = "@"
$ws.Range("A203").Value = "WC49 P5H"
$ws.Range("B203").Value = "Screw K30 no lo detecta puesto"
$ws.Range("C203").Value = "2024-06-11"
$ws.Range("D203").Value = "09:29:17"
$ws.Range("E203").Value = "Mañana"
$ws.Range("F203").Value = "09:29:17"
$ws.Range("G203").Value = "0:00:00"
$ws.Range("H203").Value = "0.08 minutos"
$rng.ClearFormats()

# Row 204
$rng = $ws.Range("A204:H204")
$rng.NumberFormat = "@"
$ws.Range("A204").Value = "WC47 NACP"
$ws.Range("B204").Value = "Fallo tolva"
$ws.Range("C204").Value = "2024-06-11"
$ws.Range("D204").Value = "09:34:07"
$ws.Range("E204").Value = "Mañana"
$ws.Range("F204").Value = "09:34:08"
$ws.Range("G204").Value = "0:00:01"
$ws.Range("H204").Value = "-0.00 minutos"
$rng.ClearFormats()

# Row 205
$rng = $ws.Range("A205:H205")
$rng.NumberFormat = "@"
$ws.Range("A205").Value = "WC49 P5H"
$ws.Range("B205").Value = "Power atascado en prensa, cuesta sacar"
$ws.Range("C205").Value = "2024-06-11"
$ws.Range("D205").Value = "09:37:55"
$ws.Range("E205").Value = "Mañana"
$ws.Range("F205").Value = "09:37:56"
$ws.Range("G205").Value = "0:00:01"
$ws.Range("H205").Value = "-0.01 minutos"
$rng.ClearFormats()

# Row 206
$rng = $ws.Range("A206:H206")
$rng.NumberFormat = "@"
$ws.Range("A206").Value = "WV50 FILTER"
$ws.Range("B206").Value = "NOK Soldadura Plástico"
$ws.Range("C206").Value = "2024-06-11"
$ws.Range("D206").Value = "09:40:13"
$ws.Range("E206").Value = "Mañana"
$ws.Range("F206").Value = "09:40:16"
$ws.Range("G206").Value = "0:00:03"
$ws.Range("H206").Value = "-0.00 minutos"
$rng.ClearFormats()

# Row 207
$rng = $ws.Range("A207:H207")
$rng.NumberFormat = "@"
$ws.Range("A207").Value = "WV50 FILTER"
$ws.Range("B207").Value = "Robot no coloca bien filter en palet"
$ws.Range("C207").Value = "2024-06-11"
$ws.Range("D207").Value = "09:40:15"
$ws.Range("E207").Value = "Mañana"
$ws.Range("F207").Value = "09:40:16"
$ws.Range("G207").Value = "0:00:01"
$ws.Range("H207").Value = "0.01 minutos"
$rng.ClearFormats()

# Row 208
$rng = $ws.Range("A208:H208")
$rng.NumberFormat = "@"
$ws.Range("A208").Value = "WV50 FILTER"
$ws.Range("B208").Value = "No coloca bien la pcb"
$ws.Range("C208").Value = "2024-06-11"
$ws.Range("D208").Value = "09:42:19"
$ws.Range("E208").Value = "Mañana"
$ws.Range("F208").Value = "09:42:20"
$ws.Range("G208").Value = "0:00:01"
$ws.Range("H208").Value = "0.35 minutos"
$rng.ClearFormats()

# Row 209
$rng = $ws.Range("A209:H209")
$rng.NumberFormat = "@"
$ws.Range("A209").Value = "WV50 FILTER"
$ws.Range("B209").Value = "Repeat funcional"
$ws.Range("C209").Value = "2024-06-11"
$ws.Range("D209").Value = "09:42:23"
$ws.Range("E209").Value = "Mañana"
$ws.Range("F209").Value = "09:42:23"
$ws.Range("G209").Value = "0:00:00"
$ws.Range("H209").Value = "0.27 minutos"
$rng.ClearFormats()
